$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.123.93"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.932.46"
$ws.Range("E3").Value = "  -4.92%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.13"
$ws.Range("E5").Value = "  -3.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("E6").Value = "  -5.67%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.43"
$ws.Range("E8").Value = "  -12.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.361"
$ws.Range("E9").Value = "  -8.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.17"
$ws.Range("E10").Value = "  -5.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0819"
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.816"
$ws.Range("E13").Value = "  -8.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.221.89"
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.67"
$ws.Range("E15").Value = "  -12.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.15"
$ws.Range("E16").Value = "  -8.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.16"
$ws.Range("E17").Value = "  -8.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.946.79"
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.053.01"
$ws.Range("E19").Value = "  -2.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.41"
$ws.Range("E20").Value = "  -4.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.39"
$ws.Range("E22").Value = "  -4.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  -8.94%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  -4.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  -7.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.28"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.06"
$ws.Range("E29").Value = "  -6.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.115"
$ws.Range("E30").Value = "  -17.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.116"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("E32").Value = "  -6.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.63"
$ws.Range("E33").Value = "  -8.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0617"
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.24"
$ws.Range("E35").Value = "  -6.73%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.99"
$ws.Range("E37").Value = "  -9.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.12"
$ws.Range("E39").Value = "  -11.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").Value = "  -12.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0961"
$ws.Range("E41").Value = "  -5.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.87"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.16"
$ws.Range("E43").Value = "  -8.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0206"
$ws.Range("E44").Value = "  -4.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.43"
$ws.Range("E45").Value = "  -9.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.335.63"
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.02"
$ws.Range("E47").Value = "  -10.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.23"
$ws.Range("E48").Value = "  -6.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.74"
$ws.Range("E49").Value = "  -7.94%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.81"
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.72"
$ws.Range("E51").Value = "  +0.59%  "
